$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row for new "database requests" block (row 9)
$ws.Range("A9").Value = "ID"
$ws.Range("B9").Value = "username"
$ws.Range("C9").Value = "passw"
$ws.Range("D9").Value = "currency"
$ws.Range("E9").Value = "currencyValue"
$ws.Range("F9").Value = "pair"
$ws.Range("G9").Value = "pairValue"
$ws.Range("H9").Value = "tokenValue"

# Type-hint row (row 10)
$ws.Range("B10").Value = "str"
$ws.Range("C10").Value = "str"
$ws.Range("D10").Value = "str"
$ws.Range("E10").Value = "float"
$ws.Range("F10").Value = "str"
$ws.Range("G10").Value = "float"
$ws.Range("H10").Value = "float"

# Column widths that Excel auto-fit (best-fit) after the new data was entered
$ws.Columns.Item(4).ColumnWidth = 7.166666666666667
$ws.Columns.Item(5).ColumnWidth = 11.833333333333334

# Restore selection to a single cell as seen in the final file
$ws.Range("N2").Select()
